$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto values (prices & 1h volume change) as produced by the
# scheduled GitHub Actions data refresh.

$ws.Range('D2').Value = '26.602.58'
$ws.Range('E2').Value = '  +3.13%  '
$ws.Range('D3').Value = '1.858.89'
$ws.Range('E3').Value = '  +2.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '272.98'
$ws.Range('E5').Value = '  -2.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5274'
$ws.Range('E7').Value = '  +3.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3370'
$ws.Range('E8').Value = '  -5.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06793'
$ws.Range('E9').Value = '  +1.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.86'
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7939'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07734'
$ws.Range('E12').Value = '  -2.06%  '
$ws.Range('D13').Value = '1.890.68'
$ws.Range('E13').Value = '  +3.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '89.67'
$ws.Range('E14').Value = '  +1.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.131'
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9992'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.40'
$ws.Range('E17').Value = '  +1.79%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007988'
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = '26.640.11'
$ws.Range('E20').Value = '  +3.08%  '
$ws.Range('D21').Value = '2.124.72'
$ws.Range('E21').Value = '  +3.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.728'
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.982'
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.110'
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.355'
$ws.Range('E25').Value = '  +4.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '145.70'
$ws.Range('E26').Value = '  +2.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.655'
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.15'
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '111.91'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.303'
$ws.Range('E30').Value = '  -0.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.291'
$ws.Range('E31').Value = '  +1.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08887'
$ws.Range('E32').Value = '  +1.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04910'
$ws.Range('E33').Value = '  -0.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.157'
$ws.Range('E34').Value = '  +1.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7281'
$ws.Range('E35').Value = '  -0.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.878'
$ws.Range('E36').Value = '  -0.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.222'
$ws.Range('E37').Value = '  +1.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.328'
$ws.Range('E38').Value = '  -1.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01840'
$ws.Range('E39').Value = '  -1.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5088'
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9376'
$ws.Range('E41').Value = '  -3.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '116.00'
$ws.Range('E42').Value = '  +0.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.139'
$ws.Range('E43').Value = '  -1.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.995'
$ws.Range('E44').Value = '  -0.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9996'
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4412'
$ws.Range('E46').Value = '  -3.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1323'
$ws.Range('E47').Value = '  -3.59%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.293'
$ws.Range('E48').Value = '  +0.89%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.02'
$ws.Range('E49').Value = '  -1.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05932'
$ws.Range('E50').Value = '  +1.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.477'
$ws.Range('E51').Value = '  -1.54%  '
